# Update the "quiz" marksheet with corrected marks totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: number right increases from 3 to 5
$ws.Range("B11").Value = 5

# Total row: total marks increases from 66 to 110 (out of 140, was 84)
$ws.Range("B12").Value = 110

# Correct/Total display string updates to match the new total
$ws.Range("E12").Value = "110/140"
